$d = $word.ActiveDocument
$d.Content.Find.Execute("289×9=", $true, $false, $false, $false, $false, $true, 1, $false, "485×6=", 2)
$d.Content.Find.Execute("899×8=", $true, $false, $false, $false, $false, $true, 1, $false, "121×9=", 2)
$d.Content.Find.Execute("425×4=", $true, $false, $false, $false, $false, $true, 1, $false, "931×6=", 2)
$d.Content.Find.Execute("131×7=", $true, $false, $false, $false, $false, $true, 1, $false, "713×3=", 2)
$d.Content.Find.Execute("231×7=", $true, $false, $false, $false, $false, $true, 1, $false, "422×9=", 2)
$d.Content.Find.Execute("356×2=", $true, $false, $false, $false, $false, $true, 1, $false, "114×4=", 2)
$d.Content.Find.Execute("362×8=", $true, $false, $false, $false, $false, $true, 1, $false, "467×9=", 2)
$d.Content.Find.Execute("270×4=", $true, $false, $false, $false, $false, $true, 1, $false, "856×7=", 2)
$d.Content.Find.Execute("553×4=", $true, $false, $false, $false, $false, $true, 1, $false, "848×9=", 2)
$d.Content.Find.Execute("585×8=", $true, $false, $false, $false, $false, $true, 1, $false, "727×5=", 2)
$d.Content.Find.Execute("752×4=", $true, $false, $false, $false, $false, $true, 1, $false, "714×9=", 2)
$d.Content.Find.Execute("244×5=", $true, $false, $false, $false, $false, $true, 1, $false, "933×5=", 2)
$d.Content.Find.Execute("653×5=", $true, $false, $false, $false, $false, $true, 1, $false, "438×7=", 2)
$d.Content.Find.Execute("683×9=", $true, $false, $false, $false, $false, $true, 1, $false, "506×6=", 2)
$d.Content.Find.Execute("451×9=", $true, $false, $false, $false, $false, $true, 1, $false, "236×2=", 2)
$d.Content.Find.Execute("819×4=", $true, $false, $false, $false, $false, $true, 1, $false, "784×2=", 2)
$d.Content.Find.Execute("161×5=", $true, $false, $false, $false, $false, $true, 1, $false, "932×6=", 2)
$d.Content.Find.Execute("187×5=", $true, $false, $false, $false, $false, $true, 1, $false, "102×2=", 2)
$d.Content.Find.Execute("776×3=", $true, $false, $false, $false, $false, $true, 1, $false, "856×2=", 2)
$d.Content.Find.Execute("621×8=", $true, $false, $false, $false, $false, $true, 1, $false, "651×6=", 2)
$d.Content.Find.Execute("519×6=", $true, $false, $false, $false, $false, $true, 1, $false, "573×8=", 2)
$d.Content.Find.Execute("842×5=", $true, $false, $false, $false, $false, $true, 1, $false, "538×2=", 2)
$d.Content.Find.Execute("568×8=", $true, $false, $false, $false, $false, $true, 1, $false, "975×8=", 2)
$d.Content.Find.Execute("889×4=", $true, $false, $false, $false, $false, $true, 1, $false, "404×7=", 2)
$d.Content.Find.Execute("254×9=", $true, $false, $false, $false, $false, $true, 1, $false, "318×2=", 2)
